$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 27 (shifts existing "Service Area" rows down by one)
$ws.Rows.Item(27).Insert()

# Populate the new row with the new Pantry item, keeping it in the sorted
# (alphabetical-by-Item) block together with the other Pantry rows.
$ws.Range("A27").Value = "Pantry"
$ws.Range("B27").Value = "Toilet Paper"

# Extend the worksheet's remembered sort range (B2:B26 -> B2:B27) to include
# the newly-inserted row, matching the expanded Pantry block.
$sort = $ws.Sort
$sort.SetRange($ws.Range("B2:B27"))
$sort.Apply()
